# UrbanGulal_Daily_2026-01-13.xlsx update
# A new order (#6, Sagar Borse) came in at 2026-01-13 18:47 and needs to be
# inserted as the newest row at the top of the order list (row 2, just below
# the header). All existing order rows shift down by one. The Summary sheet
# totals (Total Orders / New / Total Revenue) are refreshed to match.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Daily Orders")

# Push existing data rows (old rows 2-6) down to rows 3-7, leaving a blank
# row 2 for the new order.
$ws.Rows.Item(2).Insert()

# Fill in the new order in row 2.
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "2026-01-13 18:47"
$ws.Range("C2").Value = "Sagar Borse"
$ws.Range("D2").Value = "'7588930329"
$ws.Range("E2").Value = "Test,"
$ws.Range("F2").Value = "Girl Haldi Kunku Set x1"
$ws.Range("G2").Value = 25
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
$ws.Range("J2").Value = "'"
$ws.Range("K2").Value = "'"
$ws.Range("L2").Value = "'"

# The leading "'" above forces text storage (phone number / blank notes
# columns) without Excel auto-converting them to numbers; flip the style
# back to Normal so it doesn't leave a stray quote-prefix format behind.
$ws.Range("D2").Style = "Normal"
$ws.Range("J2").Style = "Normal"
$ws.Range("K2").Style = "Normal"
$ws.Range("L2").Style = "Normal"

# Refresh the Summary sheet roll-up: one more order, one more NEW order, and
# the new order's total revenue (25) added in.
$sw = $wb.Worksheets.Item("Summary")
$sw.Range("A2").Value = 6
$sw.Range("B2").Value = 5
$sw.Range("G2").Value = 25
